$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 7485.9614
$ws.Range("I15").Value = 7485.9614
$ws.Range("K15").Value = 22457.8842
$ws.Range("M15").Value = -22288.8842
# Row 17
$ws.Range("H17").Value = 11113487
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 11113487
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 33340461
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -33340797
# Row 41
$ws.Range("H41").Value = 1588
$ws.Range("I41").Value = 1588
$ws.Range("K41").Value = 1588
$ws.Range("M41").Value = -1148
# Row 62
$ws.Range("H62").Value = 2818.4
$ws.Range("I62").Value = 2273
$ws.Range("K62").Value = 2273
$ws.Range("M62").Value = -1649
# Row 65
$ws.Range("H65").Value = 2818.4
$ws.Range("I65").Value = 2273
$ws.Range("K65").Value = 11365
$ws.Range("M65").Value = -8245
# Row 70
$ws.Range("H70").Value = 3892
$ws.Range("I70").Value = 1657.2
$ws.Range("J70").Value = 4399.909
$ws.Range("K70").Value = 4971.6
$ws.Range("L70").Value = 13199.727
$ws.Range("M70").Value = -4701.6
$ws.Range("N70").Value = -13739.727
# Row 73
$ws.Range("H73").Value = 3892
$ws.Range("I73").Value = 1657.2
$ws.Range("J73").Value = 4399.909
$ws.Range("K73").Value = 4971.6
$ws.Range("L73").Value = 13199.727
$ws.Range("M73").Value = -4035.6
$ws.Range("N73").Value = -15071.727
# Row 138
$ws.Range("H138").Value = 1780.9744
$ws.Range("J138").Value = 3016.125
$ws.Range("L138").Value = 9048.375
$ws.Range("N138").Value = -19328.375

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1505.2616
$ws.Range("I32").Value = 1532.6666
$ws.Range("K32").Value = 1532.6666
$ws.Range("M32").Value = -1245.6666
# Row 63
$ws.Range("H63").Value = 4000
$ws.Range("I63").Value = 4000
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 4000
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -3314
$ws.Range("N63").ClearContents()
# Row 66
$ws.Range("H66").Value = 4000
$ws.Range("I66").Value = 4000
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 20000
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -16568
$ws.Range("N66").ClearContents()
# Row 88
$ws.Range("H88").Value = 2118
$ws.Range("J88").Value = 2063.4285
$ws.Range("L88").Value = 2063.4285
$ws.Range("N88").Value = -2875.4285
# Row 91
$ws.Range("H91").Value = 2118
$ws.Range("J91").Value = 2063.4285
$ws.Range("L91").Value = 2063.4285
$ws.Range("N91").Value = -4871.4285
# Row 122
$ws.Range("H122").Value = 3006.5
$ws.Range("I122").Value = 2771.5293
$ws.Range("J122").Value = 4338
$ws.Range("K122").Value = 8314.5879
$ws.Range("L122").Value = 13014
$ws.Range("M122").Value = -5864.5879
$ws.Range("N122").Value = -17914
# Row 132
$ws.Range("H132").Value = 3217.88
$ws.Range("I132").Value = 2653.9023
$ws.Range("K132").Value = 7961.706900000001
$ws.Range("M132").Value = -5431.706900000001

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 26
$ws.Range("H26").Value = 16944.5
$ws.Range("I26").Value = 16944.5
$ws.Range("K26").Value = 16944.5
$ws.Range("M26").Value = -16652.5
# Row 80
$ws.Range("H80").Value = 565.8182
$ws.Range("I80").Value = 360.25
$ws.Range("J80").Value = 683.2857
$ws.Range("K80").Value = 360.25
$ws.Range("L80").Value = 683.2857
$ws.Range("M80").Value = 637.75
$ws.Range("N80").Value = -2679.2857
# Row 83
$ws.Range("H83").Value = 565.8182
$ws.Range("I83").Value = 360.25
$ws.Range("J83").Value = 683.2857
$ws.Range("K83").Value = 1801.25
$ws.Range("L83").Value = 3416.4285
$ws.Range("M83").Value = 3190.75
$ws.Range("N83").Value = -13400.4285
# Row 86
$ws.Range("H86").Value = 1490
$ws.Range("I86").Value = 1500.6666
$ws.Range("J86").Value = 1479.3334
$ws.Range("K86").Value = 1500.6666
$ws.Range("L86").Value = 1479.3334
$ws.Range("M86").Value = -377.6666
$ws.Range("N86").Value = -3725.3334
# Row 89
$ws.Range("H89").Value = 1490
$ws.Range("I89").Value = 1500.6666
$ws.Range("J89").Value = 1479.3334
$ws.Range("K89").Value = 7503.333000000001
$ws.Range("L89").Value = 7396.666999999999
$ws.Range("M89").Value = -1887.333000000001
$ws.Range("N89").Value = -18628.667
# Row 94
$ws.Range("H94").Value = 1513.1428
$ws.Range("I94").Value = 1305.6364
$ws.Range("K94").Value = 1305.6364
$ws.Range("M94").Value = -854.6364000000001
# Row 107
$ws.Range("H107").Value = 1741.85
$ws.Range("J107").Value = 3850
$ws.Range("L107").Value = 3850
$ws.Range("N107").Value = -7690

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 5808
$ws.Range("J16").Value = 7623.75
$ws.Range("L16").Value = 7623.75
$ws.Range("N16").Value = -8197.75
# Row 22
$ws.Range("H22").Value = 367.1
$ws.Range("I22").Value = 364.44446
$ws.Range("K22").Value = 364.44446
$ws.Range("M22").Value = -14.44445999999999
# Row 58
$ws.Range("H58").Value = 7374.1
$ws.Range("I58").Value = 4421.6924
$ws.Range("J58").Value = 12857.143
$ws.Range("K58").Value = 4421.6924
$ws.Range("L58").Value = 12857.143
$ws.Range("M58").Value = -4218.6924
$ws.Range("N58").Value = -13263.143
# Row 107
$ws.Range("H107").Value = 658.7895
$ws.Range("J107").Value = 818
$ws.Range("L107").Value = 818
$ws.Range("N107").Value = -4658
# Row 113
$ws.Range("H113").Value = 5808
$ws.Range("J113").Value = 7623.75
$ws.Range("L113").Value = 7623.75
$ws.Range("N113").Value = -11963.75
# Row 129
$ws.Range("H129").Value = 69999
$ws.Range("J129").Value = 69999
$ws.Range("L129").Value = 69999
$ws.Range("N129").Value = -79999
# Row 134
$ws.Range("H134").Value = 4493.4287
$ws.Range("J134").Value = 4933.3
$ws.Range("L134").Value = 14799.9
$ws.Range("N134").Value = -19869.9
# Row 136
$ws.Range("H136").Value = 7374.1
$ws.Range("I136").Value = 4421.6924
$ws.Range("J136").Value = 12857.143
$ws.Range("K136").Value = 13265.0772
$ws.Range("L136").Value = 38571.429
$ws.Range("M136").Value = -10715.0772
$ws.Range("N136").Value = -43671.429

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 86
$ws.Range("H86").Value = 700
$ws.Range("I86").Value = 500
$ws.Range("K86").Value = 1500
$ws.Range("M86").Value = -314
# Row 89
$ws.Range("H89").Value = 700
$ws.Range("I89").Value = 500
$ws.Range("K89").Value = 4500
$ws.Range("M89").Value = 1428
# Row 107
$ws.Range("H107").Value = 1289.3334
$ws.Range("J107").Value = 1490.6
$ws.Range("L107").Value = 4471.799999999999
$ws.Range("N107").Value = -8311.799999999999
# Row 137
$ws.Range("H137").Value = 1789.125
$ws.Range("I137").Value = 1366.1428
$ws.Range("K137").Value = 4098.428400000001
$ws.Range("M137").Value = 1001.571599999999
# Row 139
$ws.Range("H139").Value = 1900.381
$ws.Range("I139").Value = 1683.8334
$ws.Range("J139").Value = 3199.6667
$ws.Range("K139").Value = 5051.5002
$ws.Range("L139").Value = 9599.000100000001
$ws.Range("M139").Value = 88.4997999999996
$ws.Range("N139").Value = -19879.0001
# Row 140
$ws.Range("H140").Value = 435067.22
$ws.Range("I140").Value = 854.619
$ws.Range("K140").Value = 2563.857
$ws.Range("M140").Value = 2616.143

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 12673.5
$ws.Range("I2").Value = 198.42857
$ws.Range("K2").Value = 198.42857
$ws.Range("M2").Value = -85.42857000000001
# Row 46
$ws.Range("H46").Value = 45458.332
$ws.Range("I46").Value = 39333.332
$ws.Range("J46").Value = 51583.332
$ws.Range("K46").Value = 39333.332
$ws.Range("L46").Value = 51583.332
$ws.Range("M46").Value = -39177.332
$ws.Range("N46").Value = -51895.332
# Row 57
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
# Row 80
$ws.Range("H80").Value = 3616.5312
$ws.Range("I80").Value = 3968.7334
$ws.Range("J80").Value = 3305.7646
$ws.Range("K80").Value = 3968.7334
$ws.Range("L80").Value = 3305.7646
$ws.Range("M80").Value = -2970.7334
$ws.Range("N80").Value = -5301.7646
# Row 83
$ws.Range("H83").Value = 3616.5312
$ws.Range("I83").Value = 3968.7334
$ws.Range("J83").Value = 3305.7646
$ws.Range("K83").Value = 19843.667
$ws.Range("L83").Value = 16528.823
$ws.Range("M83").Value = -14851.667
$ws.Range("N83").Value = -26512.823
# Row 109
$ws.Range("H109").Value = 50000
$ws.Range("J109").Value = 50000
$ws.Range("L109").Value = 50000
$ws.Range("N109").Value = -52080
# Row 122
$ws.Range("H122").Value = 3537.6875
$ws.Range("I122").Value = 3176.9333
$ws.Range("J122").Value = 8949
$ws.Range("K122").Value = 9530.7999
$ws.Range("L122").Value = 26847
$ws.Range("M122").Value = -7080.7999
$ws.Range("N122").Value = -31747

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 652
$ws.Range("J55").Value = 619.8
$ws.Range("L55").Value = 619.8
$ws.Range("N55").Value = -965.8
# Row 82
$ws.Range("H82").Value = 740
$ws.Range("I82").Value = 600
$ws.Range("J82").Value = 950
$ws.Range("K82").Value = 600
$ws.Range("L82").Value = 950
$ws.Range("M82").Value = -239
$ws.Range("N82").Value = -1672
# Row 85
$ws.Range("H85").Value = 740
$ws.Range("I85").Value = 600
$ws.Range("J85").Value = 950
$ws.Range("K85").Value = 600
$ws.Range("L85").Value = 950
$ws.Range("M85").Value = 648
$ws.Range("N85").Value = -3446
# Row 100
$ws.Range("H100").Value = 6115.8335
$ws.Range("J100").Value = 7332.3335
$ws.Range("L100").Value = 7332.3335
$ws.Range("N100").Value = -8414.333500000001

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 1651.6875
$ws.Range("I100").Value = 1624
$ws.Range("J100").Value = 1712.6
$ws.Range("K100").Value = 3248
$ws.Range("L100").Value = 3425.2
$ws.Range("M100").Value = -2707
$ws.Range("N100").Value = -4507.2
